$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the data previously on rows 15 and 16 -------------------------
# (Only columns whose values actually differ between the two rows need to
#  be touched; the rest of row 15/16 are identical so a swap is a no-op
#  for them.)

$row15_A = $ws.Cells.Item(15, 1).Value2
$row15_B = $ws.Cells.Item(15, 2).Value2
$row15_E = $ws.Cells.Item(15, 5).Value2
$row15_F = $ws.Cells.Item(15, 6).Value2
$row15_G = $ws.Cells.Item(15, 7).Value2
$row15_H = $ws.Cells.Item(15, 8).Value2
$row15_Z = $ws.Cells.Item(15, 26).Value2
$row15_AB = $ws.Cells.Item(15, 28).Value2

$row16_A = $ws.Cells.Item(16, 1).Value2
$row16_B = $ws.Cells.Item(16, 2).Value2
$row16_E = $ws.Cells.Item(16, 5).Value2
$row16_F = $ws.Cells.Item(16, 6).Value2
$row16_G = $ws.Cells.Item(16, 7).Value2
$row16_H = $ws.Cells.Item(16, 8).Value2
$row16_Z = $ws.Cells.Item(16, 26).Value2
$row16_AB = $ws.Cells.Item(16, 28).Value2

$ws.Cells.Item(15, 1).Value = $row16_A
$ws.Cells.Item(15, 2).Value = $row16_B
$ws.Cells.Item(15, 5).Value = $row16_E
$ws.Cells.Item(15, 6).Value = $row16_F
$ws.Cells.Item(15, 7).Value = $row16_G
$ws.Cells.Item(15, 8).Value = $row16_H
$ws.Cells.Item(15, 26).Value = $row16_Z
$ws.Cells.Item(15, 28).Value = $row16_AB

$ws.Cells.Item(16, 1).Value = $row15_A
$ws.Cells.Item(16, 2).Value = $row15_B
$ws.Cells.Item(16, 5).Value = $row15_E
$ws.Cells.Item(16, 6).Value = $row15_F
$ws.Cells.Item(16, 7).Value = $row15_G
$ws.Cells.Item(16, 8).Value = $row15_H
$ws.Cells.Item(16, 26).Value = $row15_Z
$ws.Cells.Item(16, 28).Value = $row15_AB

# --- Append two brand-new observation rows (18 & 19) ---------------------

$ws.Cells.Item(18, 1).Value = 112222262
$ws.Cells.Item(18, 2).Value = 85265
$ws.Cells.Item(18, 3).Value = "Ovaliderad"
$ws.Cells.Item(18, 4).Value = "LC"
$ws.Cells.Item(18, 5).Value = 1988
$ws.Cells.Item(18, 6).Value = "Kryddspindling"
$ws.Cells.Item(18, 7).Value = "Cortinarius percomis"
$ws.Cells.Item(18, 8).Value = "Fr."
$ws.Cells.Item(18, 16).Value = "Skogalund (Skogalund), Nrk"
$ws.Cells.Item(18, 17).Value = 531985
$ws.Cells.Item(18, 18).Value = 6553981
$ws.Cells.Item(18, 19).Value = 25
$ws.Cells.Item(18, 20).Value = "Örebro"
$ws.Cells.Item(18, 21).Value = "Örebro"
$ws.Cells.Item(18, 22).Value = "Närke"
$ws.Cells.Item(18, 23).Value = "Asker"
$ws.Cells.Item(18, 25).NumberFormat = "@"
$ws.Cells.Item(18, 25).Value = "2023-09-20"
$ws.Cells.Item(18, 27).NumberFormat = "@"
$ws.Cells.Item(18, 27).Value = "2023-09-20"
$ws.Cells.Item(18, 30).Value = $false
$ws.Cells.Item(18, 31).Value = $false
$ws.Cells.Item(18, 33).Value = $false
$ws.Cells.Item(18, 49).Value = "Magnus Friberg"
$ws.Cells.Item(18, 50).Value = "Magnus Friberg"

$ws.Cells.Item(19, 1).Value = 112222199
$ws.Cells.Item(19, 2).Value = 90379
$ws.Cells.Item(19, 3).Value = "Ovaliderad"
$ws.Cells.Item(19, 4).Value = "LC"
$ws.Cells.Item(19, 5).Value = 5836
$ws.Cells.Item(19, 6).Value = "Guldkremla"
$ws.Cells.Item(19, 7).Value = "Russula aurea"
$ws.Cells.Item(19, 8).Value = "Pers."
$ws.Cells.Item(19, 16).Value = "Skogalund (Skogalund), Nrk"
$ws.Cells.Item(19, 17).Value = 531985
$ws.Cells.Item(19, 18).Value = 6553981
$ws.Cells.Item(19, 19).Value = 25
$ws.Cells.Item(19, 20).Value = "Örebro"
$ws.Cells.Item(19, 21).Value = "Örebro"
$ws.Cells.Item(19, 22).Value = "Närke"
$ws.Cells.Item(19, 23).Value = "Asker"
$ws.Cells.Item(19, 25).NumberFormat = "@"
$ws.Cells.Item(19, 25).Value = "2023-09-20"
$ws.Cells.Item(19, 27).NumberFormat = "@"
$ws.Cells.Item(19, 27).Value = "2023-09-20"
$ws.Cells.Item(19, 30).Value = $false
$ws.Cells.Item(19, 31).Value = $false
$ws.Cells.Item(19, 33).Value = $false
$ws.Cells.Item(19, 49).Value = "Magnus Friberg"
$ws.Cells.Item(19, 50).Value = "Magnus Friberg"
